# refactored search tests, fixed some locators
#
# "OwnersWithSameLastname" worksheet gets a new block of test data for a
# new test case (verifyNoPaginationTest) and one extra row is inserted
# into the existing verifySearchPaginationNumbersTest block. The active
# sheet/selection in the workbook also moves onto this sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)   # "OwnersWithSameLastname"

# Make room for the new "verifyNoPaginationTest" block: insert 5 blank
# rows right before the old row 25. This pushes the old rows 25-41 down
# to 30-46 (carrying their existing formatting, e.g. the thick bottom
# border, with them).
$ws.Range("A25:A29").EntireRow.Insert()

# Insert one more blank row for the extra "Emily" owner that needs to
# land inside the (now shifted) verifySearchPaginationNumbersTest block,
# right after "Anthony"/"Monona" (old row 28, now row 33) and before
# "Amelia"/"Washington" (old row 29, now row 34).
$ws.Range("A34").EntireRow.Insert()

# --- New block: verifyNoPaginationTest (rows 25-29) ---
$ws.Cells.Item(25,1).Value = "Harry "
$ws.Cells.Item(25,2).Value = "Thompson"
$ws.Cells.Item(25,3).Value = "32-27 Rose Av."
$ws.Cells.Item(25,4).Value = "Boston"
$ws.Cells.Item(25,5).Value = "6085551534"
$ws.Cells.Item(25,6).Value = "verifyNoPaginationTest"

$ws.Cells.Item(26,1).Value = "Charlie"
$ws.Cells.Item(26,2).Value = "Thompson"
$ws.Cells.Item(26,3).Value = "13-12 Random St."
$ws.Cells.Item(26,4).Value = "Madison"
$ws.Cells.Item(26,5).Value = "6085552123"
$ws.Cells.Item(26,6).Value = "verifyNoPaginationTest"

$ws.Cells.Item(27,1).Value = "Jennifer"
$ws.Cells.Item(27,2).Value = "Thompson"
$ws.Cells.Item(27,3).Value = "Vrsovice 90"
$ws.Cells.Item(27,4).Value = "Praha"
$ws.Cells.Item(27,5).Value = "6085553262"
$ws.Cells.Item(27,6).Value = "verifyNoPaginationTest"

$ws.Cells.Item(28,1).Value = "John"
$ws.Cells.Item(28,2).Value = "Thompson"
$ws.Cells.Item(28,3).Value = "11-24 Long St."
$ws.Cells.Item(28,4).Value = "Windsor"
$ws.Cells.Item(28,5).Value = "6085553659"
$ws.Cells.Item(28,6).Value = "verifyNoPaginationTest"

$ws.Cells.Item(29,1).Value = "Anna"
$ws.Cells.Item(29,2).Value = "Thompson"
$ws.Cells.Item(29,3).Value = "21-22 Pernety St."
$ws.Cells.Item(29,4).Value = "Paris"
$ws.Cells.Item(29,5).Value = "6085552354"
$ws.Cells.Item(29,6).Value = "verifyNoPaginationTest"

# --- New row inside verifySearchPaginationNumbersTest (row 34) ---
$ws.Cells.Item(34,1).Value = "Emily"
$ws.Cells.Item(34,2).Value = "Thompson"
$ws.Cells.Item(34,3).Value = "4022 Annandale Rd"
$ws.Cells.Item(34,4).Value = "Washington"
$ws.Cells.Item(34,5).Value = "6085552585"
$ws.Cells.Item(34,6).Value = "verifySearchPaginationNumbersTest"

# Move the active tab/selection onto "OwnersWithSameLastname" (away from
# "OwnerWithPets"), landing on the newly edited data.
$ws.Activate()
$ws.Range("I34").Select()
